$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 179, pushing the existing 179..237 data down to 180..238
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new record
$ws.Range("A179").Value = 5
$ws.Range("B179").Value = "Macroferia Regional de Talca"
$ws.Range("C179").Value = "Maule"
$ws.Range("D179").Value = 44627
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 100112006
$ws.Range("G179").Value = "Repollo"
$ws.Range("H179").Value = "Crespo record"
$ws.Range("I179").Value = "Segunda"
$ws.Range("J179").Value = 2000
$ws.Range("K179").Value = 1000
$ws.Range("L179").Value = 1000
$ws.Range("M179").Value = 1000
$ws.Range("N179").Value = "`$/unidad"
$ws.Range("O179").Value = "Región del Maule"
$ws.Range("P179").Value = 1000
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = "Hortaliza"
